# Regenerate the save_data "K" column (G) values for erceg_lucas.xlsx
# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
#
# The K column (column G, header label "K") holds per-appearance strikeout
# counts. This recalculates/rewrites those values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0;
    3 = 0;
    4 = 3;
    5 = 1;
    6 = 2;
    7 = 1;
    8 = 1;
    9 = 2;
    10 = 2;
    11 = 2;
    12 = 2;
    13 = 1;
    14 = 1;
    15 = 1;
    16 = 2;
    17 = 0;
    18 = 1;
    19 = 0;
    20 = 1;
    21 = 2;
    22 = 4;
    23 = 0;
    24 = 2;
    25 = 1;
    26 = 2;
    27 = 1;
    28 = 0;
    29 = 2;
    30 = 1;
    31 = 1;
    32 = 1;
    33 = 0;
    34 = 0;
    35 = 1;
    36 = 1;
    37 = 1;
    38 = 1;
    39 = 0;
    40 = 1;
    41 = 2;
    42 = 1;
    43 = 0;
    44 = 2;
    45 = 0;
    46 = 1;
    47 = 1;
    48 = 1;
    49 = 0;
    50 = 1;
    51 = 2;
    52 = 0;
    53 = 2;
    54 = 3;
    55 = 3;
    56 = 2;
    57 = 2;
    58 = 1;
    59 = 2;
    60 = 0;
    61 = 2;
    62 = 1;
    63 = 0;
    64 = 2;
    65 = 0;
    66 = 1;
    67 = 1;
    68 = 1;
    72 = 2;
    73 = 2;
    75 = 1
}

foreach ($rowNum in $kValues.Keys) {
    $ws.Cells.Item([int]$rowNum, 7).Value = $kValues[$rowNum]
}
